# Apply the updated crypto market data (prices / 1h volume %, and the
# two name/link swaps) from the Mon Feb 19 09:13:03 UTC 2024 GitHub Actions refresh.
#
# Every data cell in this sheet is stored as TEXT (not a number), even when
# its contents look numeric (e.g. "351.14"). Writing straight to .Value would
# let Excel's type inference silently coerce numeric-looking strings into
# real numbers (dropping formatting like trailing zeros). To avoid that we:
#   1. force the cell to a text NumberFormat ("@") before writing,
#   2. assign the literal string to .Value,
#   3. restore the cell style to "Normal" so no visible formatting changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}
# Row 2
Set-TextValue $ws.Range("D2") "52.144.56"
Set-TextValue $ws.Range("E2") "  +0.52%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.907.14"
Set-TextValue $ws.Range("E3") "  +3.31%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.09%  "

# Row 5
Set-TextValue $ws.Range("D5") "351.14"
Set-TextValue $ws.Range("E5") "  -1.72%  "

# Row 6
Set-TextValue $ws.Range("D6") "112.13"
Set-TextValue $ws.Range("E6") "  +2.54%  "

# Row 7
Set-TextValue $ws.Range("E7") "  +0.19%  "

# Row 8
Set-TextValue $ws.Range("E8") "  +0.06%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.627"
Set-TextValue $ws.Range("E9") "  -1.19%  "

# Row 10
Set-TextValue $ws.Range("D10") "39.83"
Set-TextValue $ws.Range("E10") "  -0.61%  "

# Row 11
Set-TextValue $ws.Range("E11") "  +2.79%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.135"
Set-TextValue $ws.Range("E12") "  +0.03%  "

# Row 13
Set-TextValue $ws.Range("D13") "19.88"
Set-TextValue $ws.Range("E13") "  -0.63%  "

# Row 14
Set-TextValue $ws.Range("D14") "7.78"
Set-TextValue $ws.Range("E14") "  -0.41%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.367.45"
Set-TextValue $ws.Range("E15") "  +3.42%  "

# Row 16
Set-TextValue $ws.Range("B16") "Polygon"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D16") "0.994"
Set-TextValue $ws.Range("E16") "  +5.30%  "

# Row 17
Set-TextValue $ws.Range("B17") "WrappedEther"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D17") "2.899.61"
Set-TextValue $ws.Range("E17") "  +3.38%  "

# Row 18
Set-TextValue $ws.Range("D18") "52.201.61"
Set-TextValue $ws.Range("E18") "  +0.63%  "

# Row 19
Set-TextValue $ws.Range("E19") "  -1.12%  "

# Row 20
Set-TextValue $ws.Range("D20") "3.31"
Set-TextValue $ws.Range("E20") "  +4.46%  "

# Row 21
Set-TextValue $ws.Range("D21") "14.17"
Set-TextValue $ws.Range("E21") "  +3.31%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.0₃0977"
Set-TextValue $ws.Range("E22") "  -0.36%  "

# Row 23
Set-TextValue $ws.Range("D23") "70.66"
Set-TextValue $ws.Range("E23") "  +0.22%  "

# Row 24
Set-TextValue $ws.Range("D24") "269.39"
Set-TextValue $ws.Range("E24") "  +0.16%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.80"
Set-TextValue $ws.Range("E25") "  +1.23%  "

# Row 26
Set-TextValue $ws.Range("D26") "26.69"
Set-TextValue $ws.Range("E26") "  +1.80%  "

# Row 27
Set-TextValue $ws.Range("D27") "1.00"

# Row 28
Set-TextValue $ws.Range("D28") "0.165"
Set-TextValue $ws.Range("E28") "  +0.48%  "

# Row 29
Set-TextValue $ws.Range("D29") "10.58"
Set-TextValue $ws.Range("E29") "  +1.74%  "

# Row 30
Set-TextValue $ws.Range("D30") "37.34"
Set-TextValue $ws.Range("E30") "  -1.79%  "

# Row 31
Set-TextValue $ws.Range("D31") "2.24"
Set-TextValue $ws.Range("E31") "  +11.63%  "

# Row 32
Set-TextValue $ws.Range("D32") "6.45"
Set-TextValue $ws.Range("E32") "  +4.38%  "

# Row 33
Set-TextValue $ws.Range("E33") "  +6.72%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.0959"
Set-TextValue $ws.Range("E34") "  +10.47%  "

# Row 35
Set-TextValue $ws.Range("D35") "53.02"
Set-TextValue $ws.Range("E35") "  +2.01%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.0450"
Set-TextValue $ws.Range("E36") "  +0.85%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.999"
Set-TextValue $ws.Range("E37") "  -0.11%  "

# Row 38
Set-TextValue $ws.Range("B38") "LidoDAOToken"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D38") "3.29"
Set-TextValue $ws.Range("E38") "  +4.37%  "

# Row 39
Set-TextValue $ws.Range("D39") "18.68"
Set-TextValue $ws.Range("E39") "  -1.55%  "

# Row 40
Set-TextValue $ws.Range("B40") "ARBITRUM"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D40") "2.06"
Set-TextValue $ws.Range("E40") "  +2.59%  "

# Row 41
Set-TextValue $ws.Range("B41") "Stacks"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D41") "2.83"
Set-TextValue $ws.Range("E41") "  +12.86%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.116"
Set-TextValue $ws.Range("E42") "  +0.76%  "

# Row 43
Set-TextValue $ws.Range("D43") "23.50"
Set-TextValue $ws.Range("E43") "  +7.23%  "

# Row 44
Set-TextValue $ws.Range("D44") "2.62"
Set-TextValue $ws.Range("E44") "  +6.15%  "

# Row 45
Set-TextValue $ws.Range("D45") "120.68"
Set-TextValue $ws.Range("E45") "  +1.14%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.19"
Set-TextValue $ws.Range("E46") "  -0.45%  "

# Row 47
Set-TextValue $ws.Range("D47") "2.192.80"
Set-TextValue $ws.Range("E47") "  +3.94%  "

# Row 48
Set-TextValue $ws.Range("E48") "  +3.75%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.264"
Set-TextValue $ws.Range("E49") "  +22.90%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.962"
Set-TextValue $ws.Range("E50") "  +3.64%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.0333"
Set-TextValue $ws.Range("E51") "  +9.99%  "
